$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.003208871385164791; C = 0.00006240767534437808; D = 0.7527432677738641; E = 0.4942365360607697; G = 1.250251082895143 }
    3  = @{ B = 1.455362044514542;    C = 1.655778082260271;     D = 0.7527432677738641; E = 0.4942365360607697; G = 4.358119930609447 }
    4  = @{ B = 3.286832544864788;    C = 1.655778082260271;     D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    5  = @{ B = 0.6606524410359556;   C = 1.655778082260271;     D = 0.1494219747398047; E = 0.4942365360607697; G = 2.960089034096801 }
    6  = @{ B = 3.286832544864788;    C = 1.655778082260271;     D = 22.3905356188092;   E = 0.4942365360607697; G = 27.82738278199502 }
    7  = @{ B = 0.01293466051926884;  C = 1.655778082260271;     D = 0.7527432677738641; E = 0.4942365360607697; G = 2.915692546614173 }
    8  = @{ B = 1.455362044514542;    C = 1.655778082260271;     D = 3.537761648806719;  E = 0.4942365360607697; G = 7.143138311642302 }
    9  = @{ B = 3.286832544864788;    C = 1.655778082260271;     D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    10 = @{ B = 3.286832544864788;    C = 1.655778082260271;     D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    11 = @{ B = 3.286832544864788;    C = 1.655778082260271;     D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
